$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.736.96"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "3.383.89"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'568.52"
$ws.Range("E5").Value = "  -1.94%  "
$ws.Range("D6").Value = "'140.62"
$ws.Range("E6").Value = "  -2.91%  "
$ws.Range("D8").Value = "3.383.58"
$ws.Range("E8").Value = "  -2.04%  "
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D12").Value = "'0.398"
$ws.Range("E12").Value = "  +2.29%  "
$ws.Range("D13").Value = "3.961.92"
$ws.Range("E13").Value = "  -1.98%  "
$ws.Range("D14").Value = "'28.46"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000170"
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.385.12"
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("D18").Value = "60.824.88"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "'13.98"
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("D21").Value = "'8.97"
$ws.Range("E21").Value = "  -6.28%  "
$ws.Range("D22").Value = "'383.27"
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").Value = "'73.61"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("E26").Value = "  -6.12%  "
$ws.Range("D27").Value = "3.521.25"
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").Value = "'7.43"
$ws.Range("E30").Value = "  -2.41%  "
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.43"
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'2.14"
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "'23.63"
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("D36").Value = "'6.94"
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("D37").Value = "'166.70"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "3.413.61"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("D39").Value = "'4.99"
$ws.Range("E39").Value = "  -2.63%  "
$ws.Range("E40").Value = "  -4.74%  "
$ws.Range("D41").Value = "'27.85"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").Value = "'0.0775"
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.780"
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("D45").Value = "'41.87"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("D48").Value = "2.514.44"
$ws.Range("E48").Value = "  -2.70%  "
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("D50").Value = "'23.51"
$ws.Range("E50").Value = "  +2.24%  "
$ws.Range("E51").Value = "  -1.57%  "
